# Re-applies the latest cryptos.xlsx data refresh (GitHub Actions scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the sheet's existing inline-string
# cells) even when the value looks like a number (e.g. "98.09"), without
# leaving a numeric-format style attached to the cell afterwards.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "43.816.11"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "2.293.15"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("E4").Value = "  -0.01%  "

Set-TextValue "D5" "98.09"
$ws.Range("E5").Value = "  +4.22%  "

Set-TextValue "D6" "269.99"
$ws.Range("E6").Value = "  +0.28%  "

Set-TextValue "D7" "0.626"
$ws.Range("E7").Value = "  +1.08%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -1.36%  "

Set-TextValue "D10" "45.43"
$ws.Range("E10").Value = "  +1.16%  "

Set-TextValue "D11" "0.0937"
$ws.Range("E11").Value = "  -0.58%  "

Set-TextValue "D12" "7.89"
$ws.Range("E12").Value = "  -1.97%  "

Set-TextValue "D13" "0.106"
$ws.Range("E13").Value = "  +1.23%  "

Set-TextValue "D14" "15.85"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").Value = "2.638.73"
$ws.Range("E15").Value = "  -0.72%  "

Set-TextValue "D16" "0.857"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").Value = "2.289.03"
$ws.Range("E17").Value = "  -1.32%  "

$ws.Range("D18").Value = "43.805.49"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("E19").Value = "  +1.05%  "

Set-TextValue "D20" "6.22"
$ws.Range("E20").Value = "  -2.07%  "

Set-TextValue "D21" "72.31"
$ws.Range("E21").Value = "  +0.44%  "

Set-TextValue "D22" "2.45"
$ws.Range("E22").Value = "  +8.68%  "

Set-TextValue "D23" "233.53"
$ws.Range("E23").Value = "  -2.78%  "

Set-TextValue "D24" "9.18"
$ws.Range("E24").Value = "  -1.87%  "

Set-TextValue "D25" "2.74"
$ws.Range("E25").Value = "  +9.11%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("E29").Value = "  -0.06%  "

Set-TextValue "D30" "38.20"
$ws.Range("E30").Value = "  +0.11%  "

Set-TextValue "D31" "176.27"
$ws.Range("E31").Value = "  +2.33%  "

$ws.Range("E32").Value = "  -2.76%  "

$ws.Range("E33").Value = "  -0.24%  "

Set-TextValue "D34" "5.46"
$ws.Range("E34").Value = "  +0.35%  "

$ws.Range("E35").Value = "  +1.32%  "

Set-TextValue "D36" "4.71"
$ws.Range("E36").Value = "  +7.98%  "

Set-TextValue "D37" "0.109"
$ws.Range("E37").Value = "  +2.46%  "

Set-TextValue "D38" "0.0353"
$ws.Range("E38").Value = "  -1.54%  "

Set-TextValue "D39" "3.53"
$ws.Range("E39").Value = "  +4.82%  "

Set-TextValue "D40" "0.238"
$ws.Range("E40").Value = "  +1.53%  "

$ws.Range("E41").Value = "  +0.16%  "

Set-TextValue "D42" "1.38"
$ws.Range("E42").Value = "  +1.82%  "

Set-TextValue "D43" "12.15"
$ws.Range("E43").Value = "  +0.98%  "

Set-TextValue "D44" "64.51"
$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D45" "8.82"
$ws.Range("E45").Value = "  -3.99%  "

$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D46" "5.23"
$ws.Range("E46").Value = "  -1.92%  "

$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("E48").Value = "  +1.85%  "

Set-TextValue "D49" "98.85"
$ws.Range("E49").Value = "  -1.15%  "

Set-TextValue "D50" "0.442"
$ws.Range("E50").Value = "  +6.73%  "

Set-TextValue "D51" "1.52"
$ws.Range("E51").Value = "  +11.72%  "

